$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 203, shifting existing rows 203:251 down to 204:252
$ws.Rows.Item(203).Insert()

# Populate the new row 203 with the new record (same template columns as the
# rest of this block, with the new weekly observation's own values)
$ws.Cells.Item(203, 1).Value = 3
$ws.Cells.Item(203, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(203, 3).Value = "Coquimbo"
$ws.Cells.Item(203, 4).Value = 44543
$ws.Cells.Item(203, 5).Value = 5
$ws.Cells.Item(203, 6).Value = 100112043
$ws.Cells.Item(203, 7).Value = "Pepino ensalada"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 150
$ws.Cells.Item(203, 11).Value = 6500
$ws.Cells.Item(203, 12).Value = 7000
$ws.Cells.Item(203, 13).Value = 6733
$ws.Cells.Item(203, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(203, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(203, 16).Value = 96
$ws.Cells.Item(203, 17).Value = 70
$ws.Cells.Item(203, 18).Value = "Hortaliza"
